$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before F ("Maintainability")
#    old F/G/H/I (empty / Result / empty / Notes) shift right to G/H/I/J
# ---------------------------------------------------------------------------
$ws.Columns("F:F").Insert()

# Carry over the look of the neighbouring "Boilerplate" column onto the new one
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("E2").Copy($ws.Range("F2"))
$ws.Range("E3:E9").Copy($ws.Range("F3:F9"))

# ---------------------------------------------------------------------------
# 2. Header row (row 1) values
# ---------------------------------------------------------------------------
$ws.Range("B1").Value2 = "Learning curve"
$ws.Range("C1").Value2 = "Lib-Size"
$ws.Range("D1").Value2 = "Community"
$ws.Range("E1").Value2 = "Boilerplate"
$ws.Range("F1").Value2 = "Maintainability"

# ---------------------------------------------------------------------------
# 3. Weight row (row 2) values
# ---------------------------------------------------------------------------
$ws.Range("B2").Value2 = 4
$ws.Range("C2").Value2 = 1
$ws.Range("D2").Value2 = 3
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 4

# ---------------------------------------------------------------------------
# 4. Data rows 3-9 : new "Maintainability" ratings (col F) + a couple of
#    rating tweaks that came along with the new column
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2 = "DIY RxJS`nService / Facade"
$ws.Range("F3").Value2 = 4

$ws.Range("F4").Value2 = 3

$ws.Range("F5").Value2 = 2

$ws.Range("B6").Value2 = 2
$ws.Range("D6").Value2 = 4
$ws.Range("F6").Value2 = 1

$ws.Range("B7").Value2 = 2
$ws.Range("F7").Value2 = 1

$ws.Range("D8").Value2 = 4
$ws.Range("F8").Value2 = 2

$ws.Range("F9").Value2 = 3

Write-Host "done stage 1+2+3+4"
